$d = $word.ActiveDocument

# Remember the original author so it can be restored once we are done.
$origUser = $word.UserName

# --- Locate the stable anchor offsets in the pristine paragraph, before
# any edits are made -------------------------------------------------------
# "...Please find attached:My compression project PAQJP_4A PDF summary..."
$rngProject = $d.Content
$rngProject.Find.Execute("My compression project ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterProject = $rngProject.Start + "My compression project".Length   # right after "project", before the space
$posAfterSpace   = $rngProject.End                                       # right after that space == start of "PAQJP_"

$rng4 = $d.Content
$rng4.Find.Execute("PAQJP_4", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfter4 = $rng4.End                                                   # right after "...PAQJP_4"

# --- Use tracked changes as a controlled way to add several *separate*
# runs at the same location. A plain (untracked) edit makes this host
# normalize/merge every run in the touched paragraph that shares identical
# formatting, which would collapse all of the new pieces back together
# into one run. Tracking each insertion under a different "author" keeps
# every one of them as its own <w:ins>, and once accepted they stay apart
# as individual <w:r> elements instead of being re-coalesced.
$d.TrackRevisions = $true

# Apply the right-hand edits (around "PAQJP_4") before the left-hand ones
# (around "My compression project") so the offsets computed above - taken
# from the untouched original text - stay valid for as long as they are
# still needed.

# "...PAQJP_4" -> "...PAQJP_4" + ": " + " working already " + "and it is lossless"
$word.UserName = "Editor4"
$rLossless = $d.Range($posAfter4, $posAfter4)
$rLossless.InsertAfter("and it is lossless")

$word.UserName = "Editor5"
$rWorking = $d.Range($posAfter4, $posAfter4)
$rWorking.InsertAfter(" working already ")

$word.UserName = "Editor6"
$rColon = $d.Range($posAfter4, $posAfter4)
$rColon.InsertAfter(": ")

# "My compression project " -> "My compression project" + " Quantum" + " "
$word.UserName = "Editor1"
$delSpace = $d.Range($posAfterProject, $posAfterSpace)
$delSpace.Delete()

$word.UserName = "Editor2"
$rSep = $d.Range($posAfterProject, $posAfterProject)
$rSep.InsertAfter(" ")

$word.UserName = "Editor3"
$rQuantum = $d.Range($posAfterProject, $posAfterProject)
$rQuantum.InsertAfter(" Quantum")

# Turn tracking back off and fold every tracked change into plain content.
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()
$word.UserName = $origUser

# --- Fix up character styles ----------------------------------------------
# Runs inserted right at a paragraph's trailing edge (nothing queued after
# them yet) don't pick up the surrounding character style automatically,
# so make it explicit now that the text is plain (no longer tracked).
$fixQuantum = $d.Content
$fixQuantum.Find.Execute(" Quantum", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fixQuantum.Style = "s1"

$fixLossless = $d.Content
$fixLossless.Find.Execute("and it is lossless", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fixLossless.Style = "s2"
